$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 71 (shifts the "SBS Main Indicators..." /
# hyperlink / blank rows that followed "Source:" down by one).
$ws.Rows(71).Insert()

# Remove the hyperlink that was attached to the URL cell (now at A73) and
# turn that cell into plain text living one row further down (A74), with
# a blank row left behind at A73 - matching the restructured "Source:"
# block in the target layout.
$ws.Cells.Item(73, 1).Hyperlinks.Delete()
$ws.Cells.Item(73, 1).Value = ""
$ws.Cells.Item(74, 1).Value = "http://epp.eurostat.ec.europa.eu/portal/page/portal/european_business/data/database"

# Shorten the long citation paragraphs down to the short labels that now
# sit underneath the existing "Statistics Estonia" / "SBS Eurostat" titles.
$ws.Cells.Item(78, 1).Value = "Statistics Estonia"
$ws.Cells.Item(80, 1).Value = "SBS Eurostat"
